$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.671.46"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.603.24"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D4").Value = "'0.996"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'212.32"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'29.45"
$ws.Range("E8").Value = "  +9.33%  "
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("D10").Value = "'0.0603"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").Value = "1.832.53"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "1.604.81"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "'0.554"
$ws.Range("E14").Value = "  +2.84%  "
$ws.Range("D15").Value = "29.704.89"
$ws.Range("D16").Value = "'3.79"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").Value = "'242.40"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "'8.03"
$ws.Range("E19").Value = "  +5.72%  "
$ws.Range("D20").Value = "0.0₃0699"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'4.03"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'9.52"
$ws.Range("E23").Value = "  +3.13%  "
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "'155.45"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").Value = "'15.53"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'0.0479"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").Value = "1.424.94"
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("D35").Value = "'1.57"
$ws.Range("E35").Value = "  +4.33%  "
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("E39").Value = "  +1.67%  "
$ws.Range("D40").Value = "'0.548"
$ws.Range("E40").Value = "  +2.53%  "
$ws.Range("D41").Value = "'55.39"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.0495"
$ws.Range("E42").Value = "  +5.83%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.821"
$ws.Range("E43").Value = "  +3.39%  "
$ws.Range("D44").Value = "'1.95"
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").Value = "'0.995"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +19.29%  "
$ws.Range("D47").Value = "'67.29"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("D48").Value = "'5.37"
$ws.Range("E48").Value = "  +1.49%  "
$ws.Range("D49").Value = "1.742.29"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").Value = "'86.78"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "0.0₆0103"
$ws.Range("E51").Value = "  +2.11%  "
